$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new transformer row (row 33) at the bottom of the existing table
# Columns: A=kva, B=ht_kv, C=lt_kv, D=type, E=percentage_resistance,
#          F=percentage_reactance, G=percentage_no_load_loss
$ws.Range("A33").Value = 50000
$ws.Range("B33").Value = 115
$ws.Range("C33").Value = 30
$ws.Range("D33").Value = "unclassified"
$ws.Range("E33").Value = 0.477
$ws.Range("F33").Value = 7.5
$ws.Range("G33").Value = 0.1088

# Grow Table1 so the new row becomes part of the table range
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:G33"))

# Update selection to mirror the authoring session's last active cell
$ws.Range("G34").Select()
